$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 20:19"

# --- Plain value refreshes (country order unchanged) ---
# Row 4: Estados Unidos
$ws.Range("B4").Value = 3323231
$ws.Range("C4").Value = 31445
$ws.Range("D4").Value = 1475567
$ws.Range("E4").Value = 1710591
$ws.Range("G4").Value = 402
$ws.Range("H4").Value = 137073

# Row 6: India
$ws.Range("B6").Value = 849817
$ws.Range("C6").Value = 27214
$ws.Range("D6").Value = 536223
$ws.Range("E6").Value = 290909
$ws.Range("G6").Value = 541
$ws.Range("H6").Value = 22685

# Row 12: Reino Unido
$ws.Range("B12").Value = 288953
$ws.Range("C12").Value = 820
$ws.Range("G12").Value = 148
$ws.Range("H12").Value = 44798

# Row 19: Alemania
$ws.Range("B19").Value = 199709
$ws.Range("C19").Value = 121
$ws.Range("E19").Value = 6077

# Row 31: Ecuador
$ws.Range("B31").Value = 67209
$ws.Range("C31").Value = 2191
$ws.Range("D31").Value = 30107
$ws.Range("E31").Value = 32071
$ws.Range("G31").Value = 92
$ws.Range("H31").Value = 5031

# Row 47: Israel
$ws.Range("B47").Value = 37464
$ws.Range("C47").Value = 1198
$ws.Range("D47").Value = 18814
$ws.Range("E47").Value = 18296
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 354

# Row 65: Marruecos
$ws.Range("B65").Value = 15542
$ws.Range("C65").Value = 214
$ws.Range("D65").Value = 12065
$ws.Range("E65").Value = 3232
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 245

# Row 109: Maldivas
$ws.Range("B109").Value = 2664
$ws.Range("C109").Value = 47
$ws.Range("D109").Value = 2268
$ws.Range("E109").Value = 383

# Row 110: Sri Lanka
$ws.Range("B110").Value = 2468
$ws.Range("C110").Value = 14
$ws.Range("E110").Value = 477

# Row 129: Yemen
$ws.Range("B129").Value = 1389
$ws.Range("C129").Value = 9
$ws.Range("D129").Value = 642
$ws.Range("E129").Value = 382
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 365

# Row 146: Republica del Chad
$ws.Range("D146").Value = 790
$ws.Range("E146").Value = 10

# --- Rows 124/125: Cabo Verde overtakes Sierra Leona in ranking, so the two rows swap countries ---
$ws.Range("A124").Value = "Cabo Verde"
$ws.Range("B124").Value = 1623
$ws.Range("C124").Value = 32
$ws.Range("D124").Value = 748
$ws.Range("E124").Value = 856
$ws.Range("H124").Value = 19

$ws.Range("A125").Value = "Sierra Leona"
$ws.Range("B125").Value = 1618
$ws.Range("C125").Value = 5
$ws.Range("D125").Value = 1141
$ws.Range("E125").Value = 414
$ws.Range("H125").Value = 63

# --- Rows 136-139: Montenegro jumps ahead of Mozambique/Niger/Burkina Faso ---
$ws.Range("A136").Value = "Montenegro"
$ws.Range("B136").Value = 1164
$ws.Range("C136").Value = 145
$ws.Range("D136").Value = 325
$ws.Range("E136").Value = 820
$ws.Range("H136").Value = 19

$ws.Range("A137").Value = "Mozambique"
$ws.Range("B137").Value = 1135
$ws.Range("C137").Value = 24
$ws.Range("D137").Value = 349
$ws.Range("E137").Value = 777
$ws.Range("H137").Value = 9

$ws.Range("A138").Value = "Niger"
$ws.Range("B138").Value = 1099
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 978
$ws.Range("E138").Value = 53
$ws.Range("H138").Value = 68

$ws.Range("A139").Value = "Burkina Faso"
$ws.Range("B139").Value = 1020
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 864
$ws.Range("E139").Value = 103
$ws.Range("H139").Value = 53
